$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "2025-03-15 13:19:21"
$ws.Range("E2").Value = "POST"
$ws.Range("F2").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G2").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M2").Value = 0.002

# Row 3
$ws.Range("D3").Value = "2025-03-15 13:19:21"
$ws.Range("E3").Value = "POST"
$ws.Range("F3").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G3").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M3").Value = 0.002

# Row 4
$ws.Range("D4").Value = "2025-03-15 13:19:21"
$ws.Range("E4").Value = "POST"
$ws.Range("F4").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G4").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M4").Value = 0.003
$ws.Range("N4").Value = 0
$ws.Range("Q4").Value = $true

# Row 5
$ws.Range("D5").Value = "2025-03-15 13:19:21"
$ws.Range("E5").Value = "POST"
$ws.Range("F5").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G5").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M5").Value = 0.002
$ws.Range("N5").Value = 0
$ws.Range("O5").Value = 0
$ws.Range("P5").Value = $true
$ws.Range("Q5").Value = $true

# Row 6
$ws.Range("D6").Value = "2025-03-15 13:19:21"
$ws.Range("E6").Value = "POST"
$ws.Range("F6").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G6").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M6").Value = 0.003
$ws.Range("N6").Value = 1
$ws.Range("Q6").Value = $false

# Row 7
$ws.Range("D7").Value = "2025-03-15 13:19:21"
$ws.Range("E7").Value = "POST"
$ws.Range("F7").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G7").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M7").Value = 0.002
$ws.Range("N7").Value = 1
$ws.Range("O7").Value = 1
$ws.Range("P7").Value = $false
$ws.Range("Q7").Value = $false

# Row 8
$ws.Range("D8").Value = "2025-03-15 13:19:21"
$ws.Range("E8").Value = "POST"
$ws.Range("F8").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G8").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M8").Value = 0.002

# Row 9
$ws.Range("D9").Value = "2025-03-15 13:19:21"
$ws.Range("E9").Value = "POST"
$ws.Range("F9").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G9").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M9").Value = 0.002

# Row 10
$ws.Range("D10").Value = "2025-03-15 13:19:21"
$ws.Range("E10").Value = "POST"
$ws.Range("F10").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G10").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M10").Value = 0.002

# Row 11
$ws.Range("D11").Value = "2025-03-15 13:19:21"
$ws.Range("E11").Value = "POST"
$ws.Range("F11").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G11").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M11").Value = 0.002

# Row 12
$ws.Range("D12").Value = "2025-03-15 13:19:21"
$ws.Range("E12").Value = "POST"
$ws.Range("F12").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G12").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M12").Value = 0.002

# Row 13
$ws.Range("D13").Value = "2025-03-15 13:19:21"
$ws.Range("E13").Value = "POST"
$ws.Range("F13").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G13").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M13").Value = 0.003

# Row 14
$ws.Range("D14").Value = "2025-03-15 13:19:21"
$ws.Range("E14").Value = "POST"
$ws.Range("F14").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G14").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M14").Value = 0.002

# Row 15
$ws.Range("D15").Value = "2025-03-15 13:19:21"
$ws.Range("E15").Value = "POST"
$ws.Range("F15").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G15").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M15").Value = 0.002

# Row 16
$ws.Range("D16").Value = "2025-03-15 13:19:21"
$ws.Range("E16").Value = "POST"
$ws.Range("F16").Value = "http://49.234.6.241:5230/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("G16").Value = "/memos.api.v1.UserService/CreateUserAccessToken"
$ws.Range("M16").Value = 0.002

# Row 17
$ws.Range("D17").Value = "2025-03-15 13:19:21"
$ws.Range("E17").Value = "POST"
$ws.Range("F17").Value = "http://49.234.6.241:5230/memos.api.v1.MemoService/CreateMemo"
$ws.Range("G17").Value = "/memos.api.v1.MemoService/CreateMemo"
$ws.Range("M17").Value = 0.003
